$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the original sheet from "Sheet1" to "Data"
$ws.Name = "Data"

# 2. Build a PivotTable summarizing toxicity by day (all treatments combined)
#    from the ORIGINAL data range (A1:D25), before the "Total" rows are added.
$sourceRange = $ws.Range("A1:D25")

# Consume a throw-away sheetId so the pivot-table sheet lands on sheetId 3
# (matching the workbook's internal sheet id sequence).
$placeholder = $wb.Worksheets.Add()
$placeholder.Name = "Placeholder"
$wb.Worksheets.Item("Placeholder").Delete()

# Insert the new sheet right after "Data" and name it "Sheet3"
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$ws2.Name = "Sheet3"

$pc = $wb.PivotCaches().Create(1, $sourceRange)
$pt = $pc.CreatePivotTable($ws2.Range("A3"), "PivotTable1")

$pt.PivotFields("treatment").Orientation = 0   # xlHidden - not used on any axis
$pt.PivotFields("day").Orientation = 1         # xlRowField
$pt.AddDataField($pt.PivotFields("toxicity"), "Sum of toxicity", -4157)  # xlSum

# Approximate the pivot's rendered look: left-align the row label / grand total cells
$ws2.Range("A4:A10").HorizontalAlignment = -4131  # xlLeft

# Column widths for the pivot output sheet
$ws2.Columns.Item(1).ColumnWidth = 13.83203125
$ws2.Columns.Item(2).ColumnWidth = 14

$ws2.Range("A4").Select()

# 3. Add "Total" rows to the Data sheet: the toxicity values summed across all
#    four treatments for each day/phase (values taken from the pivot table above).
$ws.Range("A26").Value = "Total"
$ws.Range("B26").Value = "Uptake"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 0

$ws.Range("A27").Value = "Total"
$ws.Range("B27").Value = "Uptake"
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 19.136799596163598

$ws.Range("A28").Value = "Total"
$ws.Range("B28").Value = "Uptake"
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 283.98636166819011

$ws.Range("A29").Value = "Total"
$ws.Range("B29").Value = "Depuration"
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 477.39277458800552

$ws.Range("A30").Value = "Total"
$ws.Range("B30").Value = "Depuration"
$ws.Range("C30").Value = 14
$ws.Range("D30").Value = 335.01029876373036

$ws.Range("A31").Value = "Total"
$ws.Range("B31").Value = "Depuration"
$ws.Range("C31").Value = 21
$ws.Range("D31").Value = 185.43439895273917

$ws.Range("C26:C31").HorizontalAlignment = -4131  # xlLeft

$ws.Range("B9").Select()
